$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the "_old" / "_new" header-suffix columns to "_FV2310" / "_FV2404"
#    (row 1 only - A1:J1 are the *_old columns, L1:U1 are the *_new columns,
#    K1 stays "diff").
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$leftCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$rightCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($leftCols[$i] + "1").Value  = $baseNames[$i] + "_FV2310"
    $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the used range A1:U60 into an Excel Table ("Table1").
#    The header row already carries its own (bold/fill/border) styling, so we
#    temporarily stash it, reset the header to the default "Normal" style
#    before creating the table (otherwise Excel bakes the pre-existing
#    formatting into a new dxf on the table's headerRowDxfId), and restore it
#    afterwards by copying the formatting back from the stash.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A62:U62")

$headerRange.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats
$headerRange.Style = "Normal"

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U60"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$scratch.Clear()

$ws.Range("A1").Select()
